$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H112").Value = 997.1111
$ws.Range("J112").Value = 997.1111
$ws.Range("L112").Value = 2991.3333
$ws.Range("N112").Value = -5207.3333

$ws.Range("H132").Value = 4915.8335
$ws.Range("I132").Value = 2899.6
$ws.Range("K132").Value = 8698.799999999999
$ws.Range("M132").Value = -6168.799999999999

$ws.Range("H135").Value = 1011.4
$ws.Range("I135").Value = 858.3333
$ws.Range("J135").Value = 1623.6666
$ws.Range("K135").Value = 7724.9997
$ws.Range("L135").Value = 14612.9994
$ws.Range("M135").Value = -5189.9997
$ws.Range("N135").Value = -19682.9994

$ws.Range("H138").Value = 7330.727
$ws.Range("J138").Value = 7993.7793
$ws.Range("L138").Value = 23981.3379
$ws.Range("N138").Value = -34261.3379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4442.6665
$ws.Range("I32").Value = 4885.727
$ws.Range("K32").Value = 4885.727
$ws.Range("M32").Value = -4598.727

$ws.Range("H132").Value = 3617.0908
$ws.Range("J132").Value = 4432.8335
$ws.Range("L132").Value = 13298.5005
$ws.Range("N132").Value = -18358.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H107").Value = 3677.2222
$ws.Range("I107").Value = 3893
$ws.Range("J107").Value = 3245.6667
$ws.Range("K107").Value = 3893
$ws.Range("L107").Value = 3245.6667
$ws.Range("M107").Value = -1973
$ws.Range("N107").Value = -7085.6667

$ws.Range("H134").Value = 5313.1816
$ws.Range("I134").Value = 5394.5
$ws.Range("K134").Value = 16183.5
$ws.Range("M134").Value = -13648.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7708.8887
$ws.Range("I16").Value = 7967.4
$ws.Range("K16").Value = 7967.4
$ws.Range("M16").Value = -7680.4

$ws.Range("H31").Value = 1455.7778
$ws.Range("I31").Value = 886
$ws.Range("J31").Value = 3450
$ws.Range("K31").Value = 886
$ws.Range("L31").Value = 3450
$ws.Range("M31").Value = -591
$ws.Range("N31").Value = -4040

$ws.Range("H34").Value = 1455.7778
$ws.Range("I34").Value = 886
$ws.Range("J34").Value = 3450
$ws.Range("K34").Value = 886
$ws.Range("L34").Value = 3450
$ws.Range("M34").Value = -684
$ws.Range("N34").Value = -3854

$ws.Range("H62").Value = 4148.8335
$ws.Range("I62").Value = 4218.6
$ws.Range("K62").Value = 4218.6
$ws.Range("M62").Value = -3594.6

$ws.Range("H65").Value = 4148.8335
$ws.Range("I65").Value = 4218.6
$ws.Range("K65").Value = 21093
$ws.Range("M65").Value = -17973

$ws.Range("H113").Value = 7708.8887
$ws.Range("I113").Value = 7967.4
$ws.Range("K113").Value = 7967.4
$ws.Range("M113").Value = -5797.4

$ws.Range("H134").Value = 5294.6924
$ws.Range("I134").Value = 5303.1816
$ws.Range("K134").Value = 15909.5448
$ws.Range("M134").Value = -13374.5448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 800
$ws.Range("I5").Value = 800
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2400
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2288
$ws.Range("N5").ClearContents()

$ws.Range("H107").Value = 554.625
$ws.Range("I107").Value = 357
$ws.Range("J107").Value = 673.2
$ws.Range("K107").Value = 1071
$ws.Range("L107").Value = 2019.6
$ws.Range("M107").Value = 849
$ws.Range("N107").Value = -5859.6

$ws.Range("H132").Value = 9499.5
$ws.Range("J132").Value = 9000
$ws.Range("L132").Value = 81000
$ws.Range("N132").Value = -86060

$ws.Range("H135").Value = 800
$ws.Range("I135").Value = 800
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7200
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4665
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7499.5
$ws.Range("I80").Value = 6798.2
$ws.Range("K80").Value = 6798.2
$ws.Range("M80").Value = -5800.2

$ws.Range("H83").Value = 7499.5
$ws.Range("I83").Value = 6798.2
$ws.Range("K83").Value = 33991
$ws.Range("M83").Value = -28999

$ws.Range("H97").Value = 2124
$ws.Range("I97").Value = 306.69232
$ws.Range("K97").Value = 306.69232
$ws.Range("M97").Value = 189.30768

$ws.Range("H107").Value = 416.33334
$ws.Range("I107").Value = 416.75
$ws.Range("J107").Value = 415.5
$ws.Range("K107").Value = 416.75
$ws.Range("L107").Value = 415.5
$ws.Range("M107").Value = 1503.25
$ws.Range("N107").Value = -4255.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2001.8
$ws.Range("I10").Value = 1752.25
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 1752.25
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -1612.25
$ws.Range("N10").Value = -3280

$ws.Range("H16").Value = 2341.0833
$ws.Range("I16").Value = 2052.4
$ws.Range("J16").Value = 3784.5
$ws.Range("K16").Value = 2052.4
$ws.Range("L16").Value = 3784.5
$ws.Range("M16").Value = -1882.4
$ws.Range("N16").Value = -4124.5

$ws.Range("H82").Value = 1385
$ws.Range("J82").Value = 1180
$ws.Range("L82").Value = 1180
$ws.Range("N82").Value = -1902

$ws.Range("H85").Value = 1385
$ws.Range("J85").Value = 1180
$ws.Range("L85").Value = 1180
$ws.Range("N85").Value = -3676

$ws.Range("H100").Value = 1300
$ws.Range("I100").Value = 1300
$ws.Range("K100").Value = 1300
$ws.Range("M100").Value = -759

$ws.Range("H122").Value = 3204
$ws.Range("I122").Value = 3204
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9612
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7162
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2819.96
$ws.Range("I132").Value = 2301.6875
$ws.Range("J132").Value = 3741.3333
$ws.Range("K132").Value = 6905.0625
$ws.Range("L132").Value = 11223.9999
$ws.Range("M132").Value = -4375.0625
$ws.Range("N132").Value = -16283.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2870.923
$ws.Range("I132").Value = 2054
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 6162
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -3632
$ws.Range("N132").Value = -15773.4284

$ws.Range("H136").Value = 12372.044
$ws.Range("I136").Value = 12494.091
$ws.Range("K136").Value = 37482.273
$ws.Range("M136").Value = -34932.273
